$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows for more binary search leetcode problems solved.
$ws.Range("A32").Value = "Leetcode"
$ws.Range("B32").Value = 33
$ws.Range("C32").Value = "Search in Rotated Sorted Array"
$ws.Range("D32").Value = "Arrays, Binary Search"
$ws.Range("E32").Value = "Medium"
$ws.Range("F32").Value = "Neetcode 150"
$ws.Range("G32").Value = "STRUGGLED"
$ws.Range("H32").Value = "16/06/2025"
$ws.Range("I32").Value = "Struggled to handle all use cases."

$ws.Range("A33").Value = "Leetcode"
$ws.Range("B33").Value = 153
$ws.Range("C33").Value = "Find Minimum in Rotated Sorted Array"
$ws.Range("D33").Value = "Arrays, Binary Search"
$ws.Range("E33").Value = "Medium"
$ws.Range("F33").Value = "Neetcode 150"
$ws.Range("G33").Value = "SOLVED"
$ws.Range("H33").Value = "16/06/2025"
$ws.Range("I33").Value = "Easy once you have the simple technique."

# Update the TOPICS column for the existing Binary Search rows (29-31)
# to include "Arrays, Binary Search" instead of just "Binary Search".
$ws.Range("D29").Value = "Arrays, Binary Search"
$ws.Range("D30").Value = "Arrays, Binary Search"
$ws.Range("D31").Value = "Arrays, Binary Search"

# Leave the selection where the author ended up after entering the new data.
$ws.Range("I34").Select()
